$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10-12: shared "=E9" formula down E10:E12, new literal dates in F10:F12 ---
$ws.Range("D10").Value2 = "CJ Coronado, Kameron Smith"
$ws.Range("E10:E12").Formula = "=E9"
$ws.Range("F10").Value2 = 45709
$ws.Range("F11").Value2 = 45711
$ws.Range("F12").Value2 = 45711

# --- Row 13: clear status (was "In Progress") ---
$ws.Range("H13").Value2 = ""

# --- Row 17: task renamed + new assignee/start date/status/yes ---
$ws.Range("C17").Value2 = "Binary Database integration with class sytem"
$ws.Range("D17").Value2 = "Robert Snyder, CJ Coronado, Ethan Bevier"
$ws.Range("E17").Formula = "=E16"
$ws.Range("H17").Value2 = "In Progress"
$ws.Range("I17").Value2 = "Yes"

# --- Row 18: new task row ---
$ws.Range("B18").Value2 = 3.2
$ws.Range("C18").Value2 = "refine Micro UI"
$ws.Range("D18").Value2 = "Donovan Ester, Kameron Smith, Peter Jiayu Zhang"
$ws.Range("E18").Formula = "=E16"
$ws.Range("H18").Value2 = "In Progress"
$ws.Range("I18").Value2 = "Yes"

# --- Row 19: new task row ---
$ws.Range("B19").Value2 = 3.3
$ws.Range("C19").Value2 = "Floater (you'll help where it's needed)"
$ws.Range("D19").Value2 = "Timothy Barton"
$ws.Range("E19").Formula = "=E16"
$ws.Range("H19").Value2 = "In Progress"
$ws.Range("I19").Value2 = "Yes"

# --- Selection moved to D17 ---
$ws.Range("D17").Select()
